$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet tab name and header label for running "through" date
$ws.Name = "Through 2022-07-10"
$ws.Range("B1").Value = "July 2022 (through July 10)"

# Apply updated/added carjacking counts for 2022-07-10 data refresh
$ws.Range("B2").Value = 4
$ws.Range("I2").Value = 6
$ws.Range("AD2").Value = 5
$ws.Range("AK2").Value = 2
$ws.Range("B3").Value = 4
$ws.Range("P3").Value = 3
$ws.Range("AR3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("W4").Value = 2
$ws.Range("AR5").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("AK7").Value = 3
$ws.Range("P8").Value = 6
$ws.Range("AY8").Value = 3
$ws.Range("B10").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("AD16").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("B26").Value = 2
$ws.Range("I26").Value = 2
$ws.Range("B27").Value = 3
$ws.Range("B29").Value = 3
$ws.Range("AY31").Value = 1
$ws.Range("I38").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("B51").Value = 2
$ws.Range("P52").Value = 5
$ws.Range("AY59").Value = 1
$ws.Range("I62").Value = 2
$ws.Range("I65").Value = 1
$ws.Range("AY71").Value = 1
$ws.Range("I85").Value = 1
$ws.Range("I96").Value = 5
